# Scheduled market-data refresh for the Adamantoise_Profits leve-profit workbook.
# Updates currentAveragePrice(NQ/HQ) and derived LevePrice/LeveProfit columns (H-N)
# per sheet/tab with freshly polled prices; some turn-ins had no market data this run
# (value reset to 0, profit cell cleared) while others gained coverage for the first time.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1764.375
$ws.Range("J17").Value = 1764.375
$ws.Range("L17").Value = 5293.125
$ws.Range("N17").Value = -5629.125
$ws.Range("H19").Value = 1164.6923
$ws.Range("I19").Value = 1284.6
$ws.Range("J19").Value = 1089.75
$ws.Range("K19").Value = 1284.6
$ws.Range("L19").Value = 1089.75
$ws.Range("M19").Value = -1109.6
$ws.Range("N19").Value = -1439.75
$ws.Range("H100").Value = 3003.5881
$ws.Range("J100").Value = 3567.3635
$ws.Range("L100").Value = 3567.3635
$ws.Range("N100").Value = -4649.363499999999
$ws.Range("H106").Value = 4446645
$ws.Range("I106").Value = 4446645
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 4446645
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -4446014
$ws.Range("N106").ClearContents()
$ws.Range("H113").Value = 3957.7222
$ws.Range("J113").Value = 4826.8887
$ws.Range("L113").Value = 4826.8887
$ws.Range("N113").Value = -11334.8887
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()
$ws.Range("H128").Value = 141905
$ws.Range("J128").Value = 141905
$ws.Range("L128").Value = 141905
$ws.Range("N128").Value = -151865
$ws.Range("H130").Value = 116996.5
$ws.Range("J130").Value = 116996.5
$ws.Range("L130").Value = 116996.5
$ws.Range("N130").Value = -127036.5
$ws.Range("H137").Value = 13891812
$ws.Range("I137").Value = 3896
$ws.Range("K137").Value = 11688
$ws.Range("M137").Value = -9138
$ws.Range("H141").Value = 4657.909
$ws.Range("I141").Value = 4657.909
$ws.Range("K141").Value = 13973.727
$ws.Range("M141").Value = -8793.726999999999
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4407.3
$ws.Range("I61").Value = 4742.5713
$ws.Range("J61").Value = 4226.769
$ws.Range("K61").Value = 4742.5713
$ws.Range("L61").Value = 4226.769
$ws.Range("M61").Value = -4530.5713
$ws.Range("N61").Value = -4650.769
$ws.Range("H74").Value = 2832.8125
$ws.Range("J74").Value = 2261.3333
$ws.Range("L74").Value = 2261.3333
$ws.Range("N74").Value = -4009.3333
$ws.Range("H77").Value = 2832.8125
$ws.Range("J77").Value = 2261.3333
$ws.Range("L77").Value = 11306.6665
$ws.Range("N77").Value = -20042.6665
$ws.Range("H102").Value = 1399.8948
$ws.Range("I102").Value = 1152.8823
$ws.Range("K102").Value = 1152.8823
$ws.Range("M102").Value = 469.1177
$ws.Range("H132").Value = 3550
$ws.Range("I132").Value = 3511
$ws.Range("K132").Value = 10533
$ws.Range("M132").Value = -8003
$ws.Range("H136").Value = 4407.3
$ws.Range("I136").Value = 4742.5713
$ws.Range("J136").Value = 4226.769
$ws.Range("K136").Value = 14227.7139
$ws.Range("L136").Value = 12680.307
$ws.Range("M136").Value = -11677.7139
$ws.Range("N136").Value = -17780.307
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()
$ws.Range("H115").Value = 120001
$ws.Range("J115").Value = 120001
$ws.Range("L115").Value = 120001
$ws.Range("M115").Value = -123135
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()
$ws.Range("H134").Value = 4170776.5
$ws.Range("I134").Value = 4170776.5
$ws.Range("K134").Value = 12512329.5
$ws.Range("M134").Value = -12509794.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5184.25
$ws.Range("I31").Value = 5237
$ws.Range("J31").Value = 5166.6665
$ws.Range("K31").Value = 5237
$ws.Range("L31").Value = 5166.6665
$ws.Range("M31").Value = -4942
$ws.Range("N31").Value = -5756.6665
$ws.Range("H34").Value = 5184.25
$ws.Range("I34").Value = 5237
$ws.Range("J34").Value = 5166.6665
$ws.Range("K34").Value = 5237
$ws.Range("L34").Value = 5166.6665
$ws.Range("M34").Value = -5035
$ws.Range("N34").Value = -5570.6665
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()
$ws.Range("H99").Value = 3360.875
$ws.Range("I99").Value = 3198.5
$ws.Range("K99").Value = 3198.5
$ws.Range("M99").Value = -1700.5
$ws.Range("H105").Value = 2241.6
$ws.Range("I105").Value = 1942.6471
$ws.Range("K105").Value = 1942.6471
$ws.Range("M105").Value = -195.6470999999999
$ws.Range("H126").Value = 3360.875
$ws.Range("I126").Value = 3198.5
$ws.Range("K126").Value = 9595.5
$ws.Range("M126").Value = -7125.5
$ws.Range("H132").Value = 1879.6511
$ws.Range("I132").Value = 1659.7838
$ws.Range("K132").Value = 4979.3514
$ws.Range("M132").Value = -2449.3514
$ws.Range("H134").Value = 999.5
$ws.Range("I134").Value = 999.5
$ws.Range("K134").Value = 2998.5
$ws.Range("M134").Value = -463.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 99975
$ws.Range("J37").Value = 99975
$ws.Range("L37").Value = 299925
$ws.Range("N37").Value = -300149
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 481.73334
$ws.Range("I2").Value = 456.2
$ws.Range("J2").Value = 532.8
$ws.Range("K2").Value = 456.2
$ws.Range("L2").Value = 532.8
$ws.Range("M2").Value = -343.2
$ws.Range("N2").Value = -758.8
$ws.Range("H122").Value = 3045.4443
$ws.Range("I122").Value = 3045.4443
$ws.Range("K122").Value = 9136.332900000001
$ws.Range("M122").Value = -6686.332900000001
$ws.Range("H132").Value = 2952.4546
$ws.Range("I132").Value = 2952.4546
$ws.Range("K132").Value = 8857.363799999999
$ws.Range("M132").Value = -6327.363799999999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2698.8948
$ws.Range("I22").Value = 2253
$ws.Range("J22").Value = 3100.2
$ws.Range("K22").Value = 2253
$ws.Range("L22").Value = 3100.2
$ws.Range("M22").Value = -1958
$ws.Range("N22").Value = -3690.2
$ws.Range("H27").Value = 2698.8948
$ws.Range("I27").Value = 2253
$ws.Range("J27").Value = 3100.2
$ws.Range("K27").Value = 2253
$ws.Range("L27").Value = 3100.2
$ws.Range("M27").Value = -2146
$ws.Range("N27").Value = -3314.2
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("H132").Value = 3495
$ws.Range("I132").Value = 3495
$ws.Range("K132").Value = 10485
$ws.Range("M132").Value = -7955
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5790
$ws.Range("H65").Value = 5790
$ws.Range("H132").Value = 3310.6553
$ws.Range("I132").Value = 3142.647
$ws.Range("K132").Value = 9427.940999999999
$ws.Range("M132").Value = -6897.940999999999
